$d = $word.ActiveDocument

# Locate the three consecutive paragraphs that need restructuring:
#  1) "We used our own assembler ... industrial standard data."  (proofErr around "behaviour" removed)
#  2) the paragraph that used to hold only the _GoBack bookmark
#  3) the paragraph that used to start "A more detailed report ..." with the broken hyperlink/typo text
$startRng = $d.Content
$startRng.Find.Execute("We used our own assembler") | Out-Null
$startPara = $startRng.Paragraphs.Item(1)

$endRng = $d.Content
$endRng.Find.Execute("A more detailed report on our extension") | Out-Null
$endPara = $endRng.Paragraphs.Item(1)

$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = '<w:p><w:r><w:t xml:space="preserve">We used our own assembler as the operating system for the Pi, by writing an entire assembly file consisting of instructions, stating the </w:t></w:r><w:r><w:t>behaviour</w:t></w:r><w:r><w:t xml:space="preserve"> of the Pi. In the assembly file, the output pin was set then cleared, and we manually written data into memory addressed and polled from it. If our design works, then hopefully these instructions would enable the GPIO pins stated in the file to be outputting serial, industrial standard data.</w:t></w:r></w:p>' + `
        '<w:p/>' + `
        '<w:p><w:pPr><w:pStyle w:val="Body"/><w:widowControl w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>A more detailed report on our extension, as well as the i</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>mplementation, can be found at:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>https://github.com/aib13/arm11_raspberry_pi_as_DMX_show_controller</w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:pStyle w:val="Body"/><w:widowControl w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>A</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">nd a demo of our final product can be foud at: </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>https://www.youtube.com/watch?v=ou4OAfCO0nw&amp;feature=youtu.be</w:t></w:r></w:p>'

$target.InsertXML($xml)
